# IotSensorData reporting template: switch flattened-field separator from a
# single underscore to a double underscore, tweak a couple of description
# strings/numbering, and widen the columns whose header text grew.

$wb = $excel.ActiveWorkbook

# --- Sheet "semantic_aspect_model_schema": header row field names ---------
$schema = $wb.Worksheets.Item("semantic_aspect_model_schema")

$headerMap = @{
    "F1" = "sensorRuntimeData[0]__sensorGeoLocation__latitude"
    "G1" = "sensorRuntimeData[0]__sensorGeoLocation__longitude"
    "H1" = "sensorRuntimeData[0]__sensorGeoLocation__altitude"
    "I1" = "sensorRuntimeData[0]__sensorGeoLocation__geoDataTimestamp"
    "J1" = "sensorRuntimeData[0]__batteryLevel"
    "K1" = "sensorRuntimeData[0]__timestamp"
    "L1" = "sensorRuntimeData[0]__sensorData[0]__sensorType"
    "M1" = "sensorRuntimeData[0]__sensorData[0]__sensorValue"
    "N1" = "sensorRuntimeData[0]__sensorData[0]__sensorUnit"
    "O1" = "sensorRuntimeData[0]__transmissionMethod"
}

foreach ($addr in $headerMap.Keys) {
    $schema.Range($addr).Value = $headerMap[$addr]
}

# Column widths grew slightly to fit the longer (double-underscore) header
# text.
$schema.Columns.Item(10).ColumnWidth = 40.8
$schema.Columns.Item(11).ColumnWidth = 37.2
$schema.Columns.Item(15).ColumnWidth = 48

# --- Sheet "description": legend numbering, field labels, field names -----
$desc = $wb.Worksheets.Item("description")

$desc.Range("A3").Value = "1. Columns highlighted in olive green are digital twin fields."

$desc.Range("B5").Value = "Digital Twin Field Name: ownerID"
$desc.Range("B6").Value = "Digital Twin Field Name: serialNumber"
$desc.Range("B7").Value = "Digital Twin Field Name: type"
$desc.Range("B8").Value = "Digital Twin Field Name: manufacturer"

$fieldNameMap = @{
    "A10" = "sensorRuntimeData[0]__sensorGeoLocation__latitude"
    "A11" = "sensorRuntimeData[0]__sensorGeoLocation__longitude"
    "A12" = "sensorRuntimeData[0]__sensorGeoLocation__altitude"
    "A13" = "sensorRuntimeData[0]__sensorGeoLocation__geoDataTimestamp"
    "A14" = "sensorRuntimeData[0]__batteryLevel"
    "A15" = "sensorRuntimeData[0]__timestamp"
    "A16" = "sensorRuntimeData[0]__sensorData[0]__sensorType"
    "A17" = "sensorRuntimeData[0]__sensorData[0]__sensorValue"
    "A18" = "sensorRuntimeData[0]__sensorData[0]__sensorUnit"
    "A19" = "sensorRuntimeData[0]__transmissionMethod"
}

foreach ($addr in $fieldNameMap.Keys) {
    $desc.Range($addr).Value = $fieldNameMap[$addr]
}
